# Group price promotion: Name / Quantity / Price table (Orange, Tomato)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Quantity"
$ws.Range("C1").Value = "Price"

$ws.Range("A2").Value = "Orange"
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 2

$ws.Range("A3").Value = "Tomato"
$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 2

# Column widths matching the authored layout (closest values the engine's
# character-width quantization can represent).
$ws.Columns.Item(1).ColumnWidth = 21.3
$ws.Columns.Item(2).ColumnWidth = 16.3
$ws.Columns.Item(3).ColumnWidth = 13.8

$ws.Range("A3").Select()
